# Generate Report for Handoff
#
# The file "d8e6a821-6006-4af1-92ac-cc147dfbf3ac.md" (row 3 on every sheet)
# moves from "In Translation" to "Ready for handoff", and its per-language
# handoff timestamp is stamped with the moment the handoff package was cut.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: Status columns for the zh-cn / de-de languages ---
$ws_overview.Range("B3").Value = "Ready for handoff"
$ws_overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$ws_zhcn.Range("B3").Value = "Ready for handoff"
$ws_zhcn.Range("D3").Value = "2016-03-11 02:20:03"

# --- de-de sheet: Status + Latest Handoff Datetime ---
$ws_dede.Range("B3").Value = "Ready for handoff"
$ws_dede.Range("D3").Value = "2016-03-11 02:20:10"
